# Add Text to Beginning or End of Values  (new "Day-15" worksheet)

$wb = $excel.ActiveWorkbook

# --- create the new worksheet at the end of the tab strip --------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Day-15"

# --- section titles ------------------------------------------------------
$ws.Range("A1").Value = "Add Text to Beginning or End of Values"
$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("A1:D1").HorizontalAlignment = -4108
$ws.Range("A1:D1").Merge()
$ws.Rows(1).RowHeight = 15

$ws.Range("A3").Value = "#1 Flash Fill"
$ws.Range("A3:B3").Font.Bold = $true
$ws.Range("A3:B3").HorizontalAlignment = -4108
$ws.Range("A3:B3").Merge()
$ws.Rows(3).RowHeight = 15

$ws.Range("A12").Value = "#2 Formula"
$ws.Range("A12:B12").Font.Bold = $true
$ws.Range("A12:B12").HorizontalAlignment = -4108
$ws.Range("A12:B12").Merge()
$ws.Rows(12).RowHeight = 15

$ws.Range("A21").Value = "#3 Custom Format"
$ws.Range("A21:B21").Font.Bold = $true
$ws.Range("A21:B21").HorizontalAlignment = -4108
$ws.Range("A21:B21").Merge()
$ws.Rows(21).RowHeight = 15

# --- #1 Flash Fill : literal, already-filled-in values ------------------
$ws.Range("C5").Value = 270
$ws.Range("C6").Value = 450
$ws.Range("C7").Value = 781
$ws.Range("C8").Value = 562
$ws.Range("C9").Value = 124
$ws.Range("C10").Value = 892

$ws.Range("D5").Value = "ID 270"
$ws.Range("D6").Value = "ID 450"
$ws.Range("D7").Value = "ID 781"
$ws.Range("D8").Value = "ID 562"
$ws.Range("D9").Value = "ID 124"
$ws.Range("D10").Value = "ID 892"

$ws.Range("F5").Value = "Tom"
$ws.Range("F6").Value = "Luna"
$ws.Range("F7").Value = "Sara"
$ws.Range("F8").Value = "Leena"
$ws.Range("F9").Value = "Arthur"
$ws.Range("F10").Value = "James"

$ws.Range("G5").Value = "ID-Tom"
$ws.Range("G6").Value = "ID-Luna"
$ws.Range("G7").Value = "ID-Sara"
$ws.Range("G8").Value = "ID-Leena"
$ws.Range("G9").Value = "ID-Arthur"
$ws.Range("G10").Value = "ID-James"

# --- #2 Formula : CONCAT style text formulas -----------------------------
$ws.Range("C14").Value = 270
$ws.Range("C15").Value = 450
$ws.Range("C16").Value = 781
$ws.Range("C17").Value = 562
$ws.Range("C18").Value = 124
$ws.Range("C19").Value = 892

$ws.Range("D14").Formula = '="ID "&C14'
$ws.Range("D15:D19").Formula = '="ID "&C15'

$ws.Range("F14").Value = "Tom"
$ws.Range("F15").Value = "Luna"
$ws.Range("F16").Value = "Sara"
$ws.Range("F17").Value = "Leena"
$ws.Range("F18").Value = "Arthur"
$ws.Range("F19").Value = "Leila"

$ws.Range("G14").Formula = '="ID-"&F14&" Sales"'
$ws.Range("G15:G19").Formula = '="ID-"&F15&" Sales"'

# --- #3 Custom Format : number-format driven prefixing -------------------
$ws.Range("C23").Value = 270
$ws.Range("C24").Value = 450
$ws.Range("C25").Value = 781
$ws.Range("C26").Value = 562
$ws.Range("C27").Value = 124
$ws.Range("C28").Value = 892
$ws.Range("C23:C28").NumberFormat = '"ID "#'

$ws.Range("F23").Value = "Tom"
$ws.Range("F24").Value = "Luna"
$ws.Range("F25").Value = "Sara"
$ws.Range("F26").Value = "Leena"
$ws.Range("F27").Value = "Arthur"
$ws.Range("F28").Value = "Leila"
$ws.Range("F23:F28").NumberFormat = '"ID "@'

$ws.Range("I23").Value = """#"" - Numbers, ""@"" - Text"

$ws.Range("D30").Value = 270
$ws.Range("D31").Value = 450
$ws.Range("D32").Value = 781
$ws.Range("D33").Value = 562
$ws.Range("D34").Value = 124
$ws.Range("D35").Value = 892
$ws.Range("D36").Value = "Tom"
$ws.Range("D37").Value = "Luna"
$ws.Range("D38").Value = "Sara"
$ws.Range("D39").Value = "Leena"
$ws.Range("D40").Value = "Arthur"
$ws.Range("D41").Value = "Leila"
$ws.Range("D30:D41").NumberFormat = '"ID "#;;;"ID "@'

[void]$ws.Range("A1:D1").Select()
$ws.Activate()
